# "Updates on the models"
$wb = $excel.ActiveWorkbook

$global    = $wb.Worksheets.Item("Global")
$europe    = $wb.Worksheets.Item("Europe")
$catalonia = $wb.Worksheets.Item("Catalonia")

# --- Europe!A15 comment: update author / source note ----------------------
$comment = $europe.Range("A15").Comment
[void]$comment.Text("Enric:`nEDGAR`n")

# --- Europe row 15: refreshed "historic co2 emissions" series -------------
$europe.Range("C15").Value  = -0.34496931405827302
$europe.Range("D15").Value  = -0.46814242262519801
$europe.Range("E15").Value  = -0.46609129011263106
$europe.Range("F15").Value  = -0.49515421907023205
$europe.Range("G15").Value  = -0.43664204247300598
$europe.Range("H15").Value  = -0.41879503143528102
$europe.Range("I15").Value  = -0.42498349624917697
$europe.Range("J15").Value  = -0.4181721856350456
$europe.Range("K15").Value  = -0.38851950639383398
$europe.Range("L15").Value  = -0.46841254940853805
$europe.Range("M15").Value  = -0.33033722222970902
$europe.Range("N15").Value  = -0.37792358760623901
$europe.Range("O15").Value  = -0.38759988127086309
$europe.Range("P15").Value  = -0.36132755199540101
$europe.Range("Q15").Value  = -0.37335515417665699
$europe.Range("R15").Value  = -0.31262835196880601
$europe.Range("S15").Value  = -0.36857917375828403
$europe.Range("T15").Value  = -0.29009013615355
$europe.Range("U15").Value  = -0.33518804786981399
$europe.Range("V15").Value  = -0.399988360089993
$europe.Range("W15").Value  = -0.31988691783600903
$europe.Range("X15").Value  = -0.31755958479959595
$europe.Range("Y15").Value  = -0.31369644998465401
$europe.Range("Z15").Value  = -0.296045027887887
$europe.Range("AA15").Value = -0.30862722177435603
$europe.Range("AB15").Value = -0.309062138068257
$europe.Range("AC15").Value = -0.29693685476449
$europe.Range("AD15").Value = -0.27505925115152602
$europe.Range("AE15").Value = -0.24473760610362399
$europe.Range("AF15").Value = -0.24656475954628296
$europe.Range("AG15").Value = -0.26476597050030798
$europe.Range("AH15").Value = -0.22245570919959301
$europe.Range("AI15").Value = -0.22419899523269299
$europe.Range("AJ15").Value = -0.21516157963982299

# --- Per-sheet selections / which tab ends up active on save --------------
[void]$global.Activate()
[void]$global.Range("I48:J48").Select()

[void]$catalonia.Activate()
[void]$catalonia.Range("D3").Select()

[void]$europe.Activate()
[void]$europe.Range("AG15").Select()
